# Update the date/time formatting placeholders from Joda-Time based
# expressions to dateTool.format based expressions, and move the active
# cell selection from D9 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6 holds the "from - to" period template string.
$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'

# A9 holds the trip start time template string.
$ws.Range("A9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.startTime, locale, timezone)}'

# C9 holds the trip end time template string.
$ws.Range("C9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.endTime, locale, timezone)}'

# Move the active selection to B2 (was D9).
$ws.Range("B2").Select()
